$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Llegada a Datos Duros (Agregar) desde FA") is now complete (Logrado = SI)
# and gets hidden by the filter, same as the other completed rows.
$ws.Range("D5").Value = "SI"
$ws.Range("A5").EntireRow.Hidden = $true

# New row 7: another pending (NO) task, same date as row 6.
$ws.Range("A7").Value = 44461
$ws.Range("A7").NumberFormat = "d-mmm"
$ws.Range("B7").Value = "Importar avatar de producto"
$ws.Range("C7").Value = "NO"

# Re-apply the autofilter over the grown range A1:D7, keeping the
# "blank" filter on column D (Logrado) that hides completed rows.
$ws.AutoFilterMode = $false
$ws.Range("A1:D7").AutoFilter(4, @(""), 7)

# Keep the hidden _xlnm._FilterDatabase name in sync with the new range.
$fd = $wb.Names.Item("_xlnm._FilterDatabase")
$fd.RefersTo = "=Hoja1!`$A`$1:`$D`$7"

# Move the active selection to the new last row, like the author did.
$ws.Range("D7").Select()
